# The sheet contains three 21-row blocks (Brasil rows 2-22, Nordeste rows
# 23-43, Sergipe rows 44-64), each a quarterly time series in columns C
# (Trimestre / date label) and D (Valor / value). The data was rolled
# forward by one quarter: every row's C/D values become the values that
# used to belong to the following row, and a brand-new data point is
# appended at the end of each block (quarter 01/04/2024).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$groups = @(
    @{Start=2;  End=22; NewDate="01/04/2024"; NewValue=50.4425227718206},
    @{Start=23; End=43; NewDate="01/04/2024"; NewValue=43.2565479151301},
    @{Start=44; End=64; NewDate="01/04/2024"; NewValue=46.69732441471572}
)

foreach ($g in $groups) {
    $start = $g.Start
    $end = $g.End

    # Snapshot the existing C (date) and D (value) columns for this block
    # before we start overwriting anything.
    $oldC = @{}
    $oldD = @{}
    for ($r = $start; $r -le $end; $r++) {
        $oldC[$r] = $ws.Cells.Item($r, 3).Value()
        $oldD[$r] = $ws.Cells.Item($r, 4).Value()
    }

    for ($r = $start; $r -le $end; $r++) {
        if ($r -lt $end) {
            $newC = $oldC[$r + 1]
            $newD = $oldD[$r + 1]
        } else {
            $newC = $g.NewDate
            $newD = $g.NewValue
        }

        $cCell = $ws.Cells.Item($r, 3)
        $cCell.NumberFormat = "@"
        $cCell.Value = $newC
        $cCell.Style = "Normal"

        $dCell = $ws.Cells.Item($r, 4)
        if ($null -eq $newD) {
            $dCell.Value = $null
        } else {
            $dCell.Value = $newD
        }
    }
}
